$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new header row (row 61) with parameter labels in columns A:L,
# mirroring the existing data columns.
$ws.Range("A61").Value = "R1P"
$ws.Range("B61").Value = "kPL"
$ws.Range("C61").Value = "FP"
$ws.Range("D61").Value = "R1Lin"
$ws.Range("E61").Value = "kLP"
$ws.Range("F61").Value = "kMCT4"
$ws.Range("G61").Value = "R1Lex"
$ws.Range("H61").Value = "FL"
$ws.Range("I61").Value = "kMCT1"
$ws.Range("J61").Value = "k"
$ws.Range("K61").Value = "theta"
$ws.Range("L61").Value = "gamma"

# Update the view to reflect scrolling down to the newly added row.
$excel.ActiveWindow.ScrollRow = 37
$ws.Range("M61").Select()
